# Update "想去人数" (F) / "最低票价" (G) figures to the latest scrape values.
# Same underlying event rows appear on both the "展览" and "全部类型" sheets;
# the "全部类型" sheet has everything shifted down by one row.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Each entry: row on 展览 sheet, F value, G value (either can be $null to skip)
$updates = @(
    @{ Row = 5;  F = 4822; G = $null },
    @{ Row = 7;  F = 43;   G = $null },
    @{ Row = 9;  F = 484;  G = $null },
    @{ Row = 13; F = 3137; G = $null },
    @{ Row = 15; F = 116;  G = $null },
    @{ Row = 16; F = 98;   G = $null },
    @{ Row = 18; F = 2469; G = $null },
    @{ Row = 19; F = 115;  G = $null },
    @{ Row = 23; F = 24;   G = 45 },
    @{ Row = 24; F = 118;  G = $null },
    @{ Row = 25; F = $null;G = 60 },
    @{ Row = 26; F = 247;  G = $null },
    @{ Row = 27; F = 44;   G = $null }
)

foreach ($u in $updates) {
    $exhibitRow = $u.Row
    $allRow = $u.Row + 1

    if ($null -ne $u.F) {
        $wsExhibit.Range("F$exhibitRow").Value = $u.F
        $wsAll.Range("F$allRow").Value = $u.F
    }
    if ($null -ne $u.G) {
        $wsExhibit.Range("G$exhibitRow").Value = $u.G
        $wsAll.Range("G$allRow").Value = $u.G
    }
}
